# Update cryptos list data (prices and 1h volume change %) per diff.
# Numeric-looking text values in column D are prefixed with a literal
# apostrophe so Excel keeps them as text (matching the original file's
# inline-string / text formatting) instead of silently converting them
# to numbers and losing formatting such as trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.239.82'
$ws.Range("E2").Value = '  +3.37%  '
$ws.Range("D3").Value = '2.369.21'
$ws.Range("E3").Value = '  +1.67%  '
$ws.Range("E4").Value = '  -0.89%  '
$ws.Range("D5").Value = '''312.20'
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("D6").Value = '''108.89'
$ws.Range("E6").Value = '  +0.94%  '
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("E8").Value = '  -0.13%  '
$ws.Range("D9").Value = '''0.610'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("D10").Value = '''41.02'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("E11").Value = '  -0.17%  '
$ws.Range("E12").Value = '  -0.53%  '
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = '''0.978'
$ws.Range("E14").Value = '  -2.28%  '
$ws.Range("D15").Value = '2.730.51'
$ws.Range("E15").Value = '  +1.60%  '
$ws.Range("D16").Value = '''15.23'
$ws.Range("E16").Value = '  -1.43%  '
$ws.Range("D17").Value = '2.373.76'
$ws.Range("E17").Value = '  +1.67%  '
$ws.Range("D18").Value = '45.248.58'
$ws.Range("E18").Value = '  +3.58%  '
$ws.Range("D19").Value = '''14.41'
$ws.Range("E19").Value = '  +10.42%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '''0.0000106'
$ws.Range("E20").Value = '  -0.68%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = '''7.20'
$ws.Range("E21").Value = '  -4.37%  '
$ws.Range("D22").Value = '''73.07'
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").Value = '''3.52'
$ws.Range("E23").Value = '  +1.54%  '
$ws.Range("D24").Value = '''258.71'
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("E25").Value = '  +1.66%  '
$ws.Range("E26").Value = '  -0.38%  '
$ws.Range("D27").Value = '''11.08'
$ws.Range("E27").Value = '  -0.25%  '
$ws.Range("D28").Value = '''7.24'
$ws.Range("E28").Value = '  -5.67%  '
$ws.Range("E29").Value = '  +0.88%  '
$ws.Range("D30").Value = '''0.0973'
$ws.Range("E30").Value = '  +9.86%  '
$ws.Range("D31").Value = '''22.36'
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("D32").Value = '''37.39'
$ws.Range("E32").Value = '  -5.21%  '
$ws.Range("D33").Value = '''168.07'
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("D34").Value = '''3.01'
$ws.Range("E34").Value = '  +5.82%  '
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("D36").Value = '''0.117'
$ws.Range("E36").Value = '  +1.83%  '
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("E38").Value = '  +5.22%  '
$ws.Range("E39").Value = '  -2.82%  '
$ws.Range("E40").Value = '  +0.31%  '
$ws.Range("E41").Value = '  +4.75%  '
$ws.Range("D42").Value = '''99.37'
$ws.Range("E42").Value = '  -4.75%  '
$ws.Range("D43").Value = '1.895.18'
$ws.Range("E43").Value = '  +14.21%  '
$ws.Range("D44").Value = '''69.67'
$ws.Range("E44").Value = '  -2.59%  '
$ws.Range("E45").Value = '  -4.29%  '
$ws.Range("D46").Value = '''12.89'
$ws.Range("E46").Value = '  -5.39%  '
$ws.Range("E47").Value = '  -0.37%  '
$ws.Range("D48").Value = '''84.17'
$ws.Range("E48").Value = '  +10.50%  '
$ws.Range("D49").Value = '''5.65'
$ws.Range("E49").Value = '  +7.46%  '
$ws.Range("D50").Value = '''9.24'
$ws.Range("E50").Value = '  +3.14%  '
$ws.Range("D51").Value = '''110.22'
$ws.Range("E51").Value = '  -3.45%  '
